# PLU.xlsx - "Update weights loop. Test run changes"
#
# The sheet computes a Piecewise Linear Unit activation (columns A-E on
# Tabelle1). Column E holds the tunable parameters:
#   E1 = m1 (slope inside the linear region)
#   E2 = m2 (slope outside the linear region) - unchanged
#   E3 = b  (boundary output value)
#   E5 = xs = E3/E1            (computed)
#   E6 = ys = E3*(1-(E2/E1))   (computed)
# Columns B and C, and both charts, derive from these via formulas, so
# changing E1/E3 and recalculating reproduces the whole sheet + both
# chart series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the two driving weights for this test run.
$ws.Range("E1").Value = 1
$ws.Range("E3").Value = 1

# Force a full recalculation so B, C, E5 and E6 (and the chart caches,
# where supported) pick up the new weights.
$excel.CalculateFullRebuild()
$excel.Calculate()

foreach ($co in $ws.ChartObjects()) {
    $co.Chart.Refresh()
}

# Restore the window to its (new) state and move the active selection
# from E6 to E2, matching where the author left off after this run.
$excel.WindowState = -4140
$wb.Windows.Item(1).WindowState = -4140
$wb.Windows.Item(1).Left = 3990
$wb.Windows.Item(1).Top = 3270

$ws.Range("E2").Select()
